# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending-period data for the two workers (DENYS ARROYO VARGAS / HERNAN AYOLA JIMENEZ),
# alternating row by row, periods 1607..1804 (22 periods x 2 workers = 44 rows: 16..59)
$doc1 = "1048608053"
$name1 = "DENYS ARROYO VARGAS"
$doc2 = "8981082"
$name2 = "HERNAN AYOLA JIMENEZ"

$periods = @("1607","1608","1609","1610","1611","1612","1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712","1801","1802","1803","1804")

$row = 16
foreach ($p in $periods) {
    $ws.Range("C$row").Value = $doc1
    $ws.Range("D$row").Value = $name1
    $ws.Range("E$row").Value = $p
    $row = $row + 1

    $ws.Range("C$row").Value = $doc2
    $ws.Range("D$row").Value = $name2
    $ws.Range("E$row").Value = $p
    $row = $row + 1
}
